# Applies the "test donation data" edit:
#  - cleans up column B's date style on rows 3 & 4 so it matches row 2
#    (re-uses the single date format instead of the redundant one)
#  - fixes D3's donation amount
#  - adds three new donation rows (5, 6, 7), modeled on rows 2, 3 and 4
#  - moves the active selection to D7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Copy-CellFormat($srcRow, $srcCol, $dstRow, $dstCol) {
    $ws.Cells.Item($srcRow, $srcCol).Copy() | Out-Null
    $ws.Cells.Item($dstRow, $dstCol).PasteSpecial($xlPasteFormats) | Out-Null
}

function Copy-RowFormat($srcRow, $dstRow, $lastCol) {
    for ($c = 1; $c -le $lastCol; $c++) {
        Copy-CellFormat $srcRow $c $dstRow $c
    }
}

$lastCol = 54  # column BB

# --- 1. Clean up the keystroke lookup: B3/B4 should use the same date
#        style as B2 (the redundant border variant goes away). ---
Copy-CellFormat 2 2 3 2
$ws.Cells.Item(3, 2).Value = 41125

Copy-CellFormat 2 2 4 2
$ws.Cells.Item(4, 2).Value = 41126

# --- 2. Fix the donation amount in D3. ---
$ws.Cells.Item(3, 4).Value = 1235

# --- 3. Add three more fields to display: rows 5, 6 and 7, built from the
#        same formatting as rows 2, 3 and 4 respectively. ---
Copy-RowFormat 2 5 $lastCol
Copy-RowFormat 3 6 $lastCol
Copy-RowFormat 4 7 $lastCol

$ws.Rows(5).RowHeight = 32.25
$ws.Rows(6).RowHeight = 32.25
$ws.Rows(7).RowHeight = 32.25

# Row 5 (same pattern as row 2)
$ws.Cells.Item(5, 1).Value = 1050225
$ws.Cells.Item(5, 2).Value = 41127
$ws.Cells.Item(5, 3).Value = "Collector2"
$ws.Cells.Item(5, 4).Value = 12341
$ws.Cells.Item(5, 5).Value = 95308
$ws.Cells.Item(5, 6).Value = "Mrs"
$ws.Cells.Item(5, 7).Value = "Jane Doe"
$ws.Cells.Item(5, 8).Formula = "=SUM(I5:AE5)"
$ws.Cells.Item(5, 10).Value = 1440
$ws.Cells.Item(5, 12).Value = 120
$ws.Cells.Item(5, 13).Value = 50
$ws.Cells.Item(5, 14).Value = 50
$ws.Cells.Item(5, 15).Value = 12
$ws.Cells.Item(5, 16).Value = "(JulyTo Jun 12 13)"
$ws.Cells.Item(5, 21).Value = 3
$ws.Cells.Item(5, 25).Value = 10
$ws.Cells.Item(5, 26).Value = "(for previous month)"

# Row 6 (same pattern as row 3)
$ws.Cells.Item(6, 1).Value = 1050226
$ws.Cells.Item(6, 2).Value = 41128
$ws.Cells.Item(6, 3).Value = "Collector2"
$ws.Cells.Item(6, 4).Value = 12342
$ws.Cells.Item(6, 6).Value = "Mr"
$ws.Cells.Item(6, 7).Value = "ABC DEF"
$ws.Cells.Item(6, 8).Formula = "=SUM(I6:AE6)"
$ws.Cells.Item(6, 9).Value = 700
$ws.Cells.Item(6, 12).Value = 50
$ws.Cells.Item(6, 13).Value = 50
$ws.Cells.Item(6, 14).Value = 50
$ws.Cells.Item(6, 15).Value = 15
$ws.Cells.Item(6, 21).Value = 3
$ws.Cells.Item(6, 22).Value = 10
$ws.Cells.Item(6, 23).Value = 90
$ws.Cells.Item(6, 24).Value = 7
$ws.Cells.Item(6, 25).Value = 15
$ws.Cells.Item(6, 26).Value = 10
$ws.Cells.Item(6, 27).Value = "(blah blah)"

# Row 7 (same pattern as row 4)
$ws.Cells.Item(7, 1).Value = 1050227
$ws.Cells.Item(7, 2).Value = 41129
$ws.Cells.Item(7, 3).Value = "Collector1"
$ws.Cells.Item(7, 7).Value = "Michael Dean"
$ws.Cells.Item(7, 8).Formula = "=SUM(I7:AE7)"
$ws.Cells.Item(7, 27).Value = 1000

# --- 4. Move the selection to D7, matching where editing left off. ---
$ws.Range("D7").Select()
